# Output file path is removed from Input sheet for NI Scenarios
#
# The "TestResultExcelFilePath" column (column H) is deleted from the
# "ProcessPayrollForNIWeekly" and "TestReports" worksheets. Removing
# these now-unused cells also drops the two now-orphaned shared
# strings (the column header text itself and the sample output file
# path value) once the workbook is saved.

$wb = $excel.ActiveWorkbook

$wsReports = $wb.Worksheets.Item("TestReports")
$wsReports.Columns.Item(8).Select() | Out-Null
$wsReports.Columns.Item(8).Delete() | Out-Null

$wsProcess = $wb.Worksheets.Item("ProcessPayrollForNIWeekly")
$wsProcess.Columns.Item(8).Select() | Out-Null
$wsProcess.Columns.Item(8).Delete() | Out-Null
